# Update "想去人数" (want-to-go count) figures to reflect the latest
# generated data snapshot (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 1074
$ws1.Range("F11").Value = 9
$ws1.Range("F13").Value = 532
$ws1.Range("F15").Value = 12509
$ws1.Range("F16").Value = 142
$ws1.Range("F17").Value = 5499

# --- Sheet "演出" (performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 1

# --- Sheet "全部类型" (all types combined) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value  = 1074
$ws4.Range("F13").Value = 9
$ws4.Range("F15").Value = 532
$ws4.Range("F17").Value = 12509
$ws4.Range("F18").Value = 1
$ws4.Range("F20").Value = 142
$ws4.Range("F21").Value = 5499
